$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata text updates ---
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Cells changing from numeric to the shared text "0" / "***.*" (copy style+value from stable template cells in row 14) ---
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("G31"))
$ws.Range("C14").Copy($ws.Range("F33"))
$ws.Range("C14").Copy($ws.Range("G33"))
$ws.Range("E14").Copy($ws.Range("H31"))
$ws.Range("E14").Copy($ws.Range("H33"))

# --- Cells changing from shared text to numeric: fix number format first (copy from stable template cells in row 14) ---
$ws.Range("C15").NumberFormat = $ws.Range("I14").NumberFormat
$ws.Range("D15").NumberFormat = $ws.Range("I14").NumberFormat
$ws.Range("E15").NumberFormat = $ws.Range("K14").NumberFormat
$ws.Range("G15").NumberFormat = $ws.Range("I14").NumberFormat
$ws.Range("H15").NumberFormat = $ws.Range("K14").NumberFormat
$ws.Range("C27").NumberFormat = $ws.Range("I14").NumberFormat

# --- Set numeric values for the crime-data grid ---
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 30
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = 42.857142857142
$ws.Range("L15").Value = 57.894736842105
$ws.Range("M15").Value = 57.894736842105
$ws.Range("N15").Value = 50
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 185
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = -11.057692307692
$ws.Range("L16").Value = -4.145077720207
$ws.Range("M16").Value = -26
$ws.Range("N16").Value = -80
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -12
$ws.Range("I17").Value = 277
$ws.Range("J17").Value = 263
$ws.Range("K17").Value = 5.323193916349
$ws.Range("L17").Value = -1.071428571428
$ws.Range("M17").Value = 37.810945273631
$ws.Range("N17").Value = -1.77304964539
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 196
$ws.Range("J18").Value = 182
$ws.Range("K18").Value = 7.692307692307
$ws.Range("L18").Value = -21.912350597609
$ws.Range("M18").Value = -52.884615384615
$ws.Range("N18").Value = -89.39967550027
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -25.490196078431
$ws.Range("I19").Value = 585
$ws.Range("J19").Value = 642
$ws.Range("K19").Value = -8.878504672897
$ws.Range("L19").Value = -2.337228714524
$ws.Range("M19").Value = 46.616541353383
$ws.Range("N19").Value = -4.411764705882
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 410
$ws.Range("J20").Value = 340
$ws.Range("K20").Value = 20.588235294117
$ws.Range("L20").Value = 41.379310344827
$ws.Range("M20").Value = 20.588235294117
$ws.Range("N20").Value = -87.183494842138
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -8
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 119
$ws.Range("H21").Value = -5.042016806722
$ws.Range("I21").Value = 1684
$ws.Range("J21").Value = 1660
$ws.Range("K21").Value = 1.44578313253
$ws.Range("L21").Value = 2.996941896024
$ws.Range("M21").Value = 3.376304481276
$ws.Range("N21").Value = -75.61187545257
$ws.Range("G22").Value = 1
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 13.414634146341
$ws.Range("I24").Value = 1281
$ws.Range("J24").Value = 1199
$ws.Range("K24").Value = 6.839032527105
$ws.Range("L24").Value = -5.878030859662
$ws.Range("M24").Value = 30.981595092024
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 53.846153846153
$ws.Range("I25").Value = 506
$ws.Range("J25").Value = 421
$ws.Range("K25").Value = 20.190023752969
$ws.Range("L25").Value = -5.772811918063
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -27.272727272727
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 567
$ws.Range("J26").Value = 450
$ws.Range("K26").Value = 26
$ws.Range("L26").Value = 13.855421686747
$ws.Range("M26").Value = -12.5
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = 12.121212121212
$ws.Range("L27").Value = 19.354838709677
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -35
$ws.Range("N29").Value = -89.285714285714
$ws.Range("N30").Value = -88

# --- Column width adjustments for columns I (9) and J (10) ---
$ws.Columns.Item(9).ColumnWidth = 5.4
$ws.Columns.Item(10).ColumnWidth = 5.4
